# Match 18 data added (Question 2 added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (which currently ends at row 57).
# Columns: A Match | B Home Team | C Away Team | D Innings | E Batting Team |
#          F Fielding Team | G Over | H Review By | I Umpire | J Umpire Abbreviation |
#          K Decision Challenged | L Original Decision | M DRS Decision | N Batter |
#          O Bowler | P Result | Q Umpires Call

$newRows = @(
    @(18, "SRH", "CSK", 2, "SRH", "CSK", 10, "CSK", "R Pandit", "RP", "Wicket", "Not Out", "Not Out", "Shahbaz Ahmed", "M Theekshana", "Unsuccessful", "Yes"),
    @(18, "SRH", "CSK", 2, "SRH", "CSK", 14, "SRH", "R Pandit", "RP", "Wicket", "Out", "Out", "AK Markram", "MM Ali", "Unsuccessful", "No"),
    @(18, "SRH", "CSK", 2, "SRH", "CSK", 16, "CSK", "R Pandit", "RP", "Wicket", "Not Out", "Out", "Shahbaz Ahmed", "MM Ali", "Successful", "No")
)

$startRow = 58
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
